$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("WV50 FILTER", "NOK Soldadura metal",           "2024-06-05", "09:34:16", "Mañana", "09:34:16", "0:00:00", "N/A"),
    @("WV50 FILTER", "No coloca bien la pcb",         "2024-06-05", "09:34:24", "Mañana", "09:34:25", "0:00:01", "N/A"),
    @("WV50 FILTER", "Robot no coloca bien ferrita",  "2024-06-05", "09:34:53", "Mañana", "09:34:53", "0:00:00", "0.15 minutos"),
    @("WV50 FILTER", "Traza",                         "2024-06-05", "09:41:53", "Mañana", "09:41:54", "0:00:01", "N/A"),
    @("WV50 FILTER", "NOK Soldadura metal",           "2024-06-05", "09:42:10", "Mañana", "09:42:11", "0:00:01", "N/A")
)

$startRow = 20
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($col = 1; $col -le $values.Count; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        $cell.NumberFormat = "@"
        $cell.Value = $values[$col - 1]
    }
}
